$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy the header cell's style (A1: bold, bordered, centered, General
# number format) onto the whole A2:A82 date column so the bespoke
# "YYYY-MM-DD HH:MM:SS" date style is no longer used by any cell.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A82").PasteSpecial(-4122) | Out-Null

# Now replace each date-serial value in column A with its "YYYYQn" quarter
# label (computed from the underlying date serial, so this does not depend
# on any pre-baked lookup table).
for ($r = 2; $r -le 82; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2

    # Convert the Excel date serial (days since 1899-12-30) to a
    # proleptic-Gregorian (y, m, d) using the civil_from_days algorithm,
    # working entirely in integer arithmetic.
    $z = $serial - 25569 + 719468
    $era = [math]::Floor($z / 146097)
    $doe = $z - $era * 146097
    $yoe = [math]::Floor(($doe - [math]::Floor($doe / 1460) + [math]::Floor($doe / 36524) - [math]::Floor($doe / 146096)) / 365)
    $y = $yoe + $era * 400
    $doy = $doe - (365 * $yoe + [math]::Floor($yoe / 4) - [math]::Floor($yoe / 100))
    $mp = [math]::Floor((5 * $doy + 2) / 153)
    if ($mp -lt 10) { $m = $mp + 3 } else { $m = $mp - 9 }
    if ($m -le 2) { $y = $y + 1 }

    $q = [math]::Floor(($m + 1) / 3)
    $cell.Value = "$y" + "Q" + "$q"
}
